# Generate Report for Handback
#
# The "bb3e502b-971c-4e2a-a231-3c46e2280515.md" row (row 7) finished its
# handback processing for both the zh-cn and de-de target languages. The
# handback came back, but its underlying source commit was stale, so the
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# / "Error Detail" columns on row 7 get filled in (same way row 6 already
# is filled in for the failed-handback case), and a new hyperlink is added
# on the "Latest Target File" cell (I7), mirroring the existing hyperlink on
# A7 (same display text, same link target - the handback/report flow links
# back to the up-to-date source markdown file).

$wb = $excel.ActiveWorkbook

# RGB(100,149,237) == 0x6495ED, re-packed as the BGR-ish value VBA/COM
# ColorIndex expects (0x00BBGGRR) for the existing "HyperLink" font color.
$hyperlinkColor = 15570276

function Set-HandbackReady {
    param(
        [string]$sheetName,
        [string]$handbackFile,
        [string]$handbackDateTime,
        [string]$errorDetail,
        [string]$hyperlinkAddress
    )

    $ws = $wb.Worksheets.Item($sheetName)

    # I7: Latest Target File - becomes a hyperlink to the (now up to date)
    # source markdown file, same display text pattern as A7.
    $i7 = $ws.Range("I7")
    $i7.Value = "bb3e502b-971c-4e2a-a231-3c46e2280515.md"
    $i7.Font.Name = "Calibri"
    $i7.Font.Size = 11
    $i7.Font.Underline = 2
    $i7.Font.Color = $hyperlinkColor

    $ws.Hyperlinks.Add($i7, $hyperlinkAddress, "", "", "bb3e502b-971c-4e2a-a231-3c46e2280515.md")

    # J7: Latest Handback File - same file as the handoff (G7).
    $ws.Range("J7").Value = $handbackFile

    # K7: Latest Handback DateTime.
    $ws.Range("K7").Value = $handbackDateTime

    # P7: Error Detail - version mismatch notice.
    $ws.Range("P7").Value = $errorDetail
}

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/98c0e243456145359b9069e90f7b2b0f38e176c9/e2e/bb3e502b-971c-4e2a-a231-3c46e2280515.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9ce4254e2422208d0601b5d25eb28bcaac9c6763/e2e/bb3e502b-971c-4e2a-a231-3c46e2280515.md."

$latestMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9ce4254e2422208d0601b5d25eb28bcaac9c6763/e2e/bb3e502b-971c-4e2a-a231-3c46e2280515.md"

Set-HandbackReady -sheetName "zh-cn" `
    -handbackFile "bb3e502b-971c-4e2a-a231-3c46e2280515.cc706761d5f7238c24aac6fde531d360a462bbad.zh-cn.xlf" `
    -handbackDateTime "2016-09-03 08:58:43" `
    -errorDetail $errorDetail `
    -hyperlinkAddress $latestMdUrl

Set-HandbackReady -sheetName "de-de" `
    -handbackFile "bb3e502b-971c-4e2a-a231-3c46e2280515.cc706761d5f7238c24aac6fde531d360a462bbad.de-de.xlf" `
    -handbackDateTime "2016-09-03 08:58:50" `
    -errorDetail $errorDetail `
    -hyperlinkAddress $latestMdUrl

Write-Host "Handback report generated for zh-cn and de-de row 7."
